$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 92, shifting existing rows 92-114 down to 93-115
$ws.Rows(92).Insert()

# Populate the newly inserted row 92 with the new data record
$ws.Range("A92").Value = 8
$ws.Range("B92").Value = "Terminal La Palmera de La Serena"
$ws.Range("C92").Value = "Coquimbo"
$ws.Range("D92").Value = 45005
$ws.Range("E92").Value = 4
$ws.Range("F92").Value = 100114007
$ws.Range("G92").Value = "Jengibre"
$ws.Range("H92").Value = "Sin especificar"
$ws.Range("I92").Value = "Primera"
$ws.Range("J92").Value = 240
$ws.Range("K92").Value = 16000
$ws.Range("L92").Value = 17000
$ws.Range("M92").Value = 16500
$ws.Range("N92").Value = '$/caja 13 kilos'
$ws.Range("O92").Value = "Perú"
$ws.Range("P92").Value = 1269
$ws.Range("Q92").Value = 13
$ws.Range("R92").Value = "Hortaliza"
